$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Physiology")
$ws.Activate()

# Fill in "Cell #" (A), "V or I clamp" (C), "Temp" (E) for rows 7-14
$ws.Range("A7").Value = 1
for ($r = 8; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 3).Value = "V"
    $ws.Cells.Item($r, 5).Value = 34
}

$ws.Range("D14").Select()
